# Swap the full data payload (columns B:AC) between the given row pairs.
# Column A (the running id) stays untouched in each row; every other
# column (B..AC) is exchanged between the two rows of each pair, exactly
# swapping the two records with each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(136, 137),
    @(144, 147),
    @(145, 146),
    @(171, 172),
    @(178, 179),
    @(213, 214)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
